$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0298 to SCD0018
$ws.Name = "SCD0018"

# Update the TC_ID cells (B2/B3) from "DGS-313" to "SCD0018-021"
$ws.Range("B2").Value = "SCD0018-021"
$ws.Range("B3").Value = "SCD0018-021"

# Widen column B so the longer TC_ID value fits (bestFit-style width)
$ws.Columns.Item(2).ColumnWidth = 11.7

# Move the active selection from D2 to B4
$ws.Range("B4").Select()
